# "Generate Report for Handback" - mark the a.md file as handed back and
# in sync on both the zh-cn and de-de localization sheets, and refresh the
# Overview status text.

$wb = $excel.ActiveWorkbook

$blue = 15570276          # RGB(100,149,237) == OOXML color FF6495ED
$singleUnderline = 2      # xlUnderlineStyleSingle

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------
# 1. Overview sheet: status text changes from "Ready for handoff" to
#    "Handed back: in sync with en-US" for both language columns, and
#    those two columns widen to fit the new, longer text.
# ---------------------------------------------------------------------
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

$overview.Columns.Item(5).ColumnWidth = 29.15
$overview.Columns.Item(6).ColumnWidth = 29.15

# ---------------------------------------------------------------------
# 2. zh-cn and de-de detail sheets: same status text change, plus fill
#    in the handback report columns (Latest Target File / Latest
#    Handback File / Latest Handback DateTime) now that the file has
#    been handed back.
# ---------------------------------------------------------------------
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"

$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("C3").Value = "Handed back: in sync with en-US"

foreach ($ws in @($zhcn, $dede)) {
    $ws.Columns.Item(3).ColumnWidth = 29.15
    $ws.Columns.Item(10).ColumnWidth = 39.15
}

$aMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c77cc458c6de30fe2b10a9ab264f7e60c3b5712c/e2e/a.md"

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# --- zh-cn : row 2 ---
$i2 = $zhcn.Range("I2")
$zhcn.Hyperlinks.Add($i2, $aMdUrl, [Type]::Missing, [Type]::Missing, "a.md") | Out-Null
$i2.Font.Underline = $singleUnderline
$i2.Font.Color = $blue
$zhcn.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-09-05 12:46:43"
$zhcn.Range("K2").NumberFormat = $dateFmt

# --- zh-cn : row 3 ---
$i3 = $zhcn.Range("I3")
$zhcn.Hyperlinks.Add($i3, $aMdUrl, [Type]::Missing, [Type]::Missing, "a.md") | Out-Null
$i3.Font.Underline = $singleUnderline
$i3.Font.Color = $blue
$zhcn.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-09-05 12:46:43"
$zhcn.Range("K3").NumberFormat = $dateFmt

# --- de-de : row 2 ---
$i2d = $dede.Range("I2")
$dede.Hyperlinks.Add($i2d, $aMdUrl, [Type]::Missing, [Type]::Missing, "a.md") | Out-Null
$i2d.Font.Underline = $singleUnderline
$i2d.Font.Color = $blue
$dede.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("K2").Value = "2016-09-05 12:46:52"
$dede.Range("K2").NumberFormat = $dateFmt

# --- de-de : row 3 ---
$i3d = $dede.Range("I3")
$dede.Hyperlinks.Add($i3d, $aMdUrl, [Type]::Missing, [Type]::Missing, "a.md") | Out-Null
$i3d.Font.Underline = $singleUnderline
$i3d.Font.Color = $blue
$dede.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("K3").Value = "2016-09-05 12:46:52"
$dede.Range("K3").NumberFormat = $dateFmt
